$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 1292
$ws.Range("F6").Value = 317
$ws.Range("F7").Value = 1127
$ws.Range("F8").Value = 430
$ws.Range("F9").Value = 6975
$ws.Range("F11").Value = 86
$ws.Range("F12").Value = 2035
$ws.Range("F13").Value = 7871
$ws.Range("F16").Value = 5456
$ws.Range("F18").Value = 2334
$ws.Range("F19").Value = 986
$ws.Range("F21").Value = 277
$ws.Range("F25").Value = 331
$ws.Range("F26").Value = 237
$ws.Range("F27").Value = 7
$ws.Range("F28").Value = 2134
$ws.Range("F29").Value = 21
$ws.Range("F30").Value = 240
$ws.Range("F31").Value = 70
$ws.Range("F32").Value = 74
$ws.Range("F33").Value = 551
$ws.Range("F36").Value = 1432
$ws.Range("F37").Value = 27
$ws.Range("F39").Value = 2179
$ws.Range("F40").Value = 2189
$ws.Range("F41").Value = 15
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 85
$ws.Range("F3").Value = 66
$ws.Range("F4").Value = 41
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 252
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 252
$ws.Range("F6").Value = 1292
$ws.Range("F7").Value = 85
$ws.Range("F9").Value = 317
$ws.Range("F10").Value = 1127
$ws.Range("F11").Value = 430
$ws.Range("F12").Value = 6975
$ws.Range("F14").Value = 86
$ws.Range("F15").Value = 2035
$ws.Range("F16").Value = 7871
$ws.Range("F19").Value = 5456
$ws.Range("F21").Value = 2334
$ws.Range("F22").Value = 986
$ws.Range("F24").Value = 277
$ws.Range("F27").Value = 66
$ws.Range("F29").Value = 41
$ws.Range("F30").Value = 331
$ws.Range("F31").Value = 237
$ws.Range("F32").Value = 7
$ws.Range("F33").Value = 2134
$ws.Range("F34").Value = 21
$ws.Range("F35").Value = 240
$ws.Range("F36").Value = 70
$ws.Range("F37").Value = 74
$ws.Range("F38").Value = 551
$ws.Range("F42").Value = 1432
$ws.Range("F43").Value = 27
$ws.Range("F45").Value = 2179
$ws.Range("F47").Value = 2189
$ws.Range("F48").Value = 15
